$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the column-A date cells down into the new rows
$ws.Range("A35").Copy()
$ws.Range("A36:A37").PasteSpecial(-4122)  # xlPasteFormats

# Row 36: new entry for 2012-11-05 (date serial 41218)
$ws.Cells.Item(36, 1).Value = 41218
$ws.Cells.Item(36, 2).Value = 1
$ws.Cells.Item(36, 4).Value = "Manual continued"

# Row 37: new entry for 2012-11-06 (date serial 41219)
$ws.Cells.Item(37, 1).Value = 41219
$ws.Cells.Item(37, 2).Value = 1.75
$ws.Cells.Item(37, 3).Value = 0.25
$ws.Cells.Item(37, 4).Value = "Manual continued, new test case tc07"

$ws.Range("C37").Select()
